$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'90.880.65"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "'3.167.75"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'215.74"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").Value = "'627.44"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  +29.81%  "
$ws.Range("E8").Value = "  +3.24%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'3.165.31"
$ws.Range("E10").Value = "  +3.92%  "
$ws.Range("D11").Value = "'0.753"
$ws.Range("E11").Value = "  +11.00%  "
$ws.Range("D12").Value = "'0.204"
$ws.Range("E12").Value = "  +8.58%  "
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("E14").Value = "  +5.82%  "
$ws.Range("D15").Value = "'35.01"
$ws.Range("E15").Value = "  +7.94%  "
$ws.Range("D16").Value = "'90.554.69"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "'3.740.34"
$ws.Range("D18").Value = "'3.141.58"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("E19").Value = "  +10.99%  "
$ws.Range("D20").Value = "'14.38"
$ws.Range("E20").Value = "  +7.05%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'466.81"
$ws.Range("E21").Value = "  +9.78%  "
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").Value = "'0.0000212"
$ws.Range("E22").Value = "  -4.01%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'9.12"
$ws.Range("E23").Value = "  +10.44%  "
$ws.Range("B24").Value = "Polkadot"
$ws.Range("C24").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D24").Value = "'5.27"
$ws.Range("E24").Value = "  +4.79%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'93.94"
$ws.Range("E25").Value = "  +11.82%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'5.88"
$ws.Range("E26").Value = "  +8.56%  "
$ws.Range("D27").Value = "'12.23"
$ws.Range("E27").Value = "  +4.90%  "
$ws.Range("D28").Value = "'3.317.58"
$ws.Range("E28").Value = "  +3.58%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'9.20"
$ws.Range("E32").Value = "  +10.01%  "
$ws.Range("D33").Value = "'27.85"
$ws.Range("E33").Value = "  +21.86%  "
$ws.Range("D34").Value = "'523.92"
$ws.Range("E34").Value = "  +3.97%  "
$ws.Range("D35").Value = "'0.184"
$ws.Range("E35").Value = "  +35.05%  "
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("E37").Value = "  +7.40%  "
$ws.Range("D38").Value = "'6.94"
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'1.31"
$ws.Range("E39").Value = "  +5.44%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.142"
$ws.Range("E40").Value = "  +7.15%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "'0.0854"
$ws.Range("E42").Value = "  +23.70%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'0.415"
$ws.Range("E44").Value = "  +14.16%  "
$ws.Range("D45").Value = "'1.97"
$ws.Range("E45").Value = "  +7.67%  "
$ws.Range("D47").Value = "'150.17"
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D48").Value = "'45.38"
$ws.Range("E48").Value = "  +4.50%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'0.686"
$ws.Range("E49").Value = "  +16.49%  "
$ws.Range("E50").Value = "  +10.98%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'4.51"
$ws.Range("E51").Value = "  +7.35%  "
